# Generate Report for Handback
# - Marks zh-cn / de-de as handed back (in sync with en-US) instead of
#   "Ready for handoff".
# - Refreshes the "Latest Handback DateTime" stamps for both languages.
# - Clears the stale "handback file is not latest" Error Detail message
#   now that the handback is in sync.
# - Widens the Status columns / narrows the Error Detail columns to fit
#   the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns -------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

$overview.Range("E1").EntireColumn.ColumnWidth = 29.166666666666668
$overview.Range("F1").EntireColumn.ColumnWidth = 29.166666666666668

# --- zh-cn detail sheet -----------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-13 01:01:16"
$zhcn.Range("P2").Value = ""

$zhcn.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$zhcn.Range("P1").EntireColumn.ColumnWidth = 12.833333333333334

# --- de-de detail sheet ------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-13 01:01:26"
$dede.Range("P2").Value = ""

$dede.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$dede.Range("P1").EntireColumn.ColumnWidth = 12.833333333333334
